# feat: esquema exemplo para planilha de importação
#
# 1. Rename the sheet Folha1 -> Folha2
# 2. Replace the "challenge 2" sample row (row 3) and its extra shared
#    strings with a single new string "Design Figma"
# 3. Re-base the cellXfs table: Excel inserted a new default+protection
#    style (index 1) ahead of the existing ones, bake that in by touching
#    Format Cells > Protection on the header row before the data edits
# 4. Delete the now-unused second data row and its hyperlink
# 5. Reset page margins / header&footer to the Excel defaults

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Folha2"

# 2. Delete the sample "Desafio 2" row (row 3) and its hyperlink.
#    Hyperlinks.Delete() clears every hyperlink on the sheet in this
#    runtime, so snapshot the surviving one and re-add it afterwards.
$ws.Hyperlinks.Delete()

$ws.Rows(3).Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "http://doity@teste.com") | Out-Null

# 3. New cell content for E2 ("Serviços" column on the remaining row)
$ws.Range("E2").Value = "Design Figma"

# 4. Explicit default protection on the header row + formatting touch-up
#    that Excel bakes into a new cellXfs entry (fontId 0 / fillId 0 /
#    borderId 0, protection hidden=0 locked=1) used across row 1 and most
#    of row 2.
$headerRange = $ws.Range("A1:E1")
$headerRange.Locked = $true
$headerRange.FormulaHidden = $false

$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Style = "Normal"
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").Style = "Normal"
$ws.Range("E1").Style = "Normal"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Style = "Normal"

# 5. Row height / column width reset
$ws.Rows("1:2").RowHeight = 15
$ws.Columns("D").AutoFit() | Out-Null

# 6. Page margins back to Excel's printer defaults
$ws.PageSetup.LeftMargin = 0.70078740157480324 * 72
$ws.PageSetup.RightMargin = 0.70078740157480324 * 72
$ws.PageSetup.TopMargin = 0.75196850393700787 * 72
$ws.PageSetup.BottomMargin = 0.75196850393700787 * 72
$ws.PageSetup.HeaderMargin = 0.29999999999999999 * 72
$ws.PageSetup.FooterMargin = 0.29999999999999999 * 72
